# Update Leave Card 1/4/2024 4:38 PM
# Adds 9 more monthly leave-credit rows (Apr 2023 - Dec 2023, + 2 blank
# trailing rows) to the "Sheet1" leave ledger, fills in the two rows that
# were already present but still blank (Mar/Apr-2023-dated rows 153-154),
# and grows Table1 / the sheet dimension to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Fill in the EARNED value for the two existing-but-blank rows ----
$ws.Range("C153").Value2 = 1.25
$ws.Range("C154").Value2 = 1.25

# --- 2. Grow Table1 down to row 163 (9 new rows) ------------------------
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A8:K163"))

# --- 3. Clone the formatting of row 154 into the freshly added rows -----
$ws.Range("A154:K154").Copy() | Out-Null
$ws.Range("A155:K163").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- 4. Populate rows 155-161 (monthly accrual continues); 162-163 stay blank
$periodRows = 155..161
foreach ($r in $periodRows) {
    $prev = $r - 1
    $ws.Range("A$r").Formula = "=EDATE(A$prev,1)"
    $ws.Range("C$r").Value2 = 1.25
    $ws.Range("G$r").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),""."",Table1[[#This Row],[EARNED]])"
    $ws.Range("G$r").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"
}

# Rows 162 & 163 remain formatted-but-empty placeholder rows (already done
# by the format-only paste above).

# --- 5. Refresh the bottom-pane selection to follow the newly entered cell
$ws.Range("E156").Select()

$wb.Save()
